# Modify wait in all pages classes and update TestData file records
# (only the TestData1.xlsx workbook edits are relevant to this script)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the test data records stored in row 8 (shared strings) ---
$ws.Range("A8").Value = "ATestBAutomationA1"
$ws.Range("B8").Value = "ATestBAutomationA1"
$ws.Range("C8").Value = "Facility_C1091a1222"
$ws.Range("D8").Value = "Facility_C1091a1222"
$ws.Range("E8").Value = "Pharmacy_C1091a1222"
$ws.Range("F8").Value = "Pharmacy_C1091a1222"
$ws.Range("H8").Value = "Alignment Project C1091a1222"

# --- Update the sheet's active selection ---
$ws.Range("F13").Select()

# --- Widen column H slightly ---
$ws.Columns.Item(8).ColumnWidth = 28.6666666666667
